$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; this shifts all existing rows (and their
# values/formatting) down by one, preserving the original data exactly as-is.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "Variable"
$ws.Range("B1").Value = "Value"
